$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new rows below the existing data, carrying the s="1" number
#     format down from row 12 so the new metric cells match the rest of the table ---
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(15).Insert()

# --- New notes text in column D for the new rows (fills shared-string slots first) ---
$ws.Cells.Item(13, 4).Value = '"", 3days means, time to eta in seconds'
$ws.Cells.Item(14, 4).Value = '60ts, ""'
$ws.Cells.Item(15, 4).Value = '72ts, isMoored'

# --- Result labels in column A (reverse row order, then fix the renamed result_2.0 label) ---
$ws.Cells.Item(15, 1).Value = 'result_2.3'
$ws.Cells.Item(14, 1).Value = 'result_2.2'
$ws.Cells.Item(13, 1).Value = 'result_2.1'
$ws.Cells.Item(12, 1).Value = 'result_2.0'

# --- git-lfs instructions in new column H ---
$ws.Cells.Item(12, 8).Value = 'sudo apt install git-lfs'
$ws.Cells.Item(13, 8).Value = 'git lfs install'

# --- Numeric metrics for the new rows ---
$ws.Cells.Item(13, 2).Value = 0.76579529521828305
$ws.Cells.Item(13, 3).Value = 102.74959

$ws.Cells.Item(14, 2).Value = 0.65849783034842602
$ws.Cells.Item(14, 3).Value = 103.86183

$ws.Cells.Item(15, 2).Value = 0.77123279476968698
$ws.Cells.Item(15, 3).Value = 101.11895

# --- Wrap text for the long notes column C on the new "ts"/"isMoored" rows ---
$ws.Range("C14:C15").WrapText = $true

# --- New column H width ---
$ws.Columns.Item(8).ColumnWidth = 22.5

# --- Final selection, matching the saved workbook state ---
$ws.Range("H13").Select()
